# Update the "dSF" (column F) values on Sheet1 to reflect the repulled
# data / mean calculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -2
    6  = 1
    7  = -2
    9  = 1
    11 = -5
    12 = 1
    14 = -3
    15 = 7
    16 = 2
    17 = 1
    22 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
